$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add two new worksheets at the end of the workbook: "same_elements" and
# "partly_same" (new array-sort benchmark scenario: arrays whose elements
# are all identical / mostly identical).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSame = $wb.Worksheets.Add($null, $lastSheet)
$wsSame.Name = "same_elements"

$wsPartly = $wb.Worksheets.Add($null, $wsSame)
$wsPartly.Name = "partly_same"

function Fill-BenchmarkSheet {
    param($ws, $rows)

    # header row: array sizes benchmarked
    $ws.Range("B1").Value = 5
    $ws.Range("C1").Value = 50
    $ws.Range("D1").Value = 500
    $ws.Range("E1").Value = 5000
    $ws.Range("F1").Value = 50000
    $ws.Range("G1").Value = 500000

    $labels = @("byte", "int", "string", "date")
    for ($i = 0; $i -lt 4; $i++) {
        $r = $i + 2
        $ws.Range("A$r").Value = $labels[$i]
        $data = $rows[$i]
        $ws.Range("B$r").Value = $data[0]
        $ws.Range("C$r").Value = $data[1]
        $ws.Range("D$r").Value = $data[2]
        $ws.Range("E$r").Value = $data[3]
        $ws.Range("F$r").Value = $data[4]
        $ws.Range("G$r").Value = $data[5]
    }
}

$sameElementsRows = @(
    @(0, 0, 0,      0.002504, 0.025525, 0.216259),
    @(0, 0, 0.0005, 0.002536, 0.022021, 0.217694),
    @(0, 0, 0,      0.002001, 0.022023, 0.217727),
    @(0, 0, 0.000499, 0.002504, 0.022023, 0.244257)
)

$partlySameRows = @(
    @(0, 0, 0.003502, 0.041543, 0.577805, 6.565403),
    @(0, 0, 0.003505, 0.048586, 0.66245,  8.123995),
    @(0, 0, 0.003503, 0.048513, 0.655852, 8.622764),
    @(0, 0, 0.003504, 0.048548, 0.640826, 9.330857)
)

Fill-BenchmarkSheet $wsSame $sameElementsRows
Fill-BenchmarkSheet $wsPartly $partlySameRows

# "same_elements" becomes the newly active/selected sheet.
$wsSame.Activate()
